$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="296.46"'
$ws.Range("E2").Formula = '="3.62%"'
$ws.Range("G2").Formula = '="16"'
$ws.Range("D3").Formula = '="41.43"'
$ws.Range("E3").Formula = '="2.86%"'
$ws.Range("G3").Formula = '="16"'
$ws.Range("D4").Formula = '="5.052"'
$ws.Range("E4").Formula = '="0.47%"'
$ws.Range("G4").Formula = '="16"'
$ws.Range("D5").Formula = '="0.07473"'
$ws.Range("E5").Formula = '="2.12%"'
$ws.Range("G5").Formula = '="16"'
$ws.Range("D6").Formula = '="1.575"'
$ws.Range("E6").Formula = '="2.14%"'
$ws.Range("G6").Formula = '="16"'
$ws.Range("D7").Formula = '="0.9288"'
$ws.Range("E7").Formula = '="1.71%"'
$ws.Range("G7").Formula = '="16"'
$ws.Range("G8").Formula = '="16"'
$ws.Range("D9").Formula = '="0.1197"'
$ws.Range("E9").Formula = '="-0.76%"'
$ws.Range("G9").Formula = '="16"'
$ws.Range("D10").Formula = '="0.1811"'
$ws.Range("E10").Formula = '="4.60%"'
$ws.Range("G10").Formula = '="16"'
$ws.Range("D11").Formula = '="0.08871"'
$ws.Range("E11").Formula = '="2.88%"'
$ws.Range("G11").Formula = '="16"'
$ws.Range("D12").Formula = '="0.04301"'
$ws.Range("E12").Formula = '="3.16%"'
$ws.Range("G12").Formula = '="16"'
$ws.Range("D13").Formula = '="0.1051"'
$ws.Range("E13").Formula = '="0.10%"'
$ws.Range("G13").Formula = '="16"'
$ws.Range("D14").Formula = '="0.001275"'
$ws.Range("E14").Formula = '="0.05%"'
$ws.Range("G14").Formula = '="16"'
$ws.Range("D15").Formula = '="0.005972"'
$ws.Range("E15").Formula = '="1.25%"'
$ws.Range("G15").Formula = '="16"'
$ws.Range("E16").Formula = '="-1.34%"'
$ws.Range("G16").Formula = '="16"'
$ws.Range("D17").Formula = '="4.372"'
$ws.Range("E17").Formula = '="2.06%"'
$ws.Range("G17").Formula = '="16"'
$ws.Range("D18").Formula = '="0.3302"'
$ws.Range("E18").Formula = '="1.09%"'
$ws.Range("G18").Formula = '="16"'
$ws.Range("D19").Formula = '="8.051"'
$ws.Range("E19").Formula = '="5.79%"'
$ws.Range("G19").Formula = '="16"'
$ws.Range("D20").Formula = '="0.1379"'
$ws.Range("E20").Formula = '="2.74%"'
$ws.Range("G20").Formula = '="16"'
$ws.Range("D21").Formula = '="0.2964"'
$ws.Range("E21").Formula = '="2.77%"'
$ws.Range("G21").Formula = '="16"'
$ws.Range("D22").Formula = '="0.04027"'
$ws.Range("E22").Formula = '="5.14%"'
$ws.Range("G22").Formula = '="16"'
$ws.Range("D23").Formula = '="0.001266"'
$ws.Range("E23").Formula = '="-0.31%"'
$ws.Range("G23").Formula = '="16"'
$ws.Range("D24").Formula = '="0.003870"'
$ws.Range("E24").Formula = '="2.17%"'
$ws.Range("G24").Formula = '="16"'
$ws.Range("D25").Formula = '="0.0001228"'
$ws.Range("E25").Formula = '="-4.24%"'
$ws.Range("G25").Formula = '="16"'
$ws.Range("D26").Formula = '="0.0003717"'
$ws.Range("E26").Formula = '="-0.35%"'
$ws.Range("G26").Formula = '="16"'
$ws.Range("G27").Formula = '="16"'
$ws.Range("G28").Formula = '="16"'
$ws.Range("G29").Formula = '="16"'
$ws.Range("G30").Formula = '="16"'
$ws.Range("G31").Formula = '="16"'
$ws.Range("G32").Formula = '="16"'
$ws.Range("G33").Formula = '="16"'
$ws.Range("G34").Formula = '="16"'
$ws.Range("G35").Formula = '="16"'
$ws.Range("G36").Formula = '="16"'
$ws.Range("G37").Formula = '="16"'
$ws.Range("D38").Formula = '="0.02407"'
$ws.Range("E38").Formula = '="3.70%"'
$ws.Range("G38").Formula = '="16"'
$ws.Range("D39").Formula = '="0.05188"'
$ws.Range("E39").Formula = '="4.05%"'
$ws.Range("G39").Formula = '="16"'
$ws.Range("D40").Formula = '="0.006691"'
$ws.Range("E40").Formula = '="30.95%"'
$ws.Range("G40").Formula = '="16"'
$ws.Range("D41").Formula = '="0.007790"'
$ws.Range("E41").Formula = '="1.32%"'
$ws.Range("G41").Formula = '="16"'
$ws.Range("D42").Formula = '="0.1322"'
$ws.Range("E42").Formula = '="4.09%"'
$ws.Range("G42").Formula = '="16"'
$ws.Range("D43").Formula = '="0.007360"'
$ws.Range("E43").Formula = '="-0.07%"'
$ws.Range("G43").Formula = '="16"'
$ws.Range("D44").Formula = '="0.007795"'
$ws.Range("E44").Formula = '="3.76%"'
$ws.Range("G44").Formula = '="16"'
$ws.Range("D45").Formula = '="0.3214"'
$ws.Range("E45").Formula = '="4.16%"'
$ws.Range("G45").Formula = '="16"'
$ws.Range("D46").Formula = '="0.00006346"'
$ws.Range("E46").Formula = '="-2.11%"'
$ws.Range("G46").Formula = '="16"'
$ws.Range("D47").Formula = '="0.00000000749"'
$ws.Range("E47").Formula = '="-0.34%"'
$ws.Range("G47").Formula = '="16"'
$ws.Range("D48").Formula = '="0.04688"'
$ws.Range("E48").Formula = '="-81.38%"'
$ws.Range("G48").Formula = '="16"'
$ws.Range("D49").Formula = '="0.004194"'
$ws.Range("E49").Formula = '="-0.18%"'
$ws.Range("G49").Formula = '="16"'
$ws.Range("D50").Formula = '="0.00002097"'
$ws.Range("E50").Formula = '="-0.34%"'
$ws.Range("G50").Formula = '="16"'
$ws.Range("D51").Formula = '="0.0001997"'
$ws.Range("E51").Formula = '="-0.34%"'
$ws.Range("G51").Formula = '="16"'

$rng = $ws.Range("D2:G51")
$rng.Copy()
$rng.PasteSpecial(-4163)
